$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellXml($cell, $segments) {
  $runsXml = ""
  $first = $true
  foreach ($seg in $segments) {
    if (-not $first) {
      $runsXml += "<w:br/>"
    }
    $first = $false
    $text = $seg
    # Mirror Word's own behaviour: only emit xml:space="preserve" when the
    # text has leading/trailing whitespace (or is empty) that must survive.
    $needsPreserve = ($text -ne $text.Trim()) -or ($text -eq "")
    if ($needsPreserve) {
      $runsXml += '<w:t xml:space="preserve">' + $text + '</w:t>'
    } else {
      $runsXml += '<w:t>' + $text + '</w:t>'
    }
  }
  $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + $runsXml + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $cell.Range.InsertXML($xml) | Out-Null
}

# Each lattice-multiplication cell holds 5 lines of text, joined by <w:br/>:
#   1) "A x B"           -- the two factors
#   2) "  d1    d2"      -- digits of the second factor spaced above the grid
#   3) "  ----"          -- the separator bar (never changes)
#   4) "d|    |"         -- first digit of the first factor
#   5) "d|    |"         -- second digit of the first factor
# The table keeps its 5 rows x 3 columns shape; only the exercise numbers
# inside each cell are refreshed to the newly generated problems.

Set-CellXml $t.Cell(1,1) @("49 x 56", "  5    6", "  ----", "4|    |", "9|    |")
Set-CellXml $t.Cell(1,2) @("97 x 87", "  8    7", "  ----", "9|    |", "7|    |")
Set-CellXml $t.Cell(1,3) @("60 x 86", "  8    6", "  ----", "6|    |", "0|    |")

Set-CellXml $t.Cell(2,1) @("37 x 20", "  2    0", "  ----", "3|    |", "7|    |")
Set-CellXml $t.Cell(2,2) @("76 x 24", "  2    4", "  ----", "7|    |", "6|    |")
Set-CellXml $t.Cell(2,3) @("11 x 94", "  9    4", "  ----", "1|    |", "1|    |")

Set-CellXml $t.Cell(3,1) @("49 x 89", "  8    9", "  ----", "4|    |", "9|    |")
Set-CellXml $t.Cell(3,2) @("34 x 82", "  8    2", "  ----", "3|    |", "4|    |")
Set-CellXml $t.Cell(3,3) @("18 x 89", "  8    9", "  ----", "1|    |", "8|    |")

Set-CellXml $t.Cell(4,1) @("99 x 35", "  3    5", "  ----", "9|    |", "9|    |")
Set-CellXml $t.Cell(4,2) @("20 x 87", "  8    7", "  ----", "2|    |", "0|    |")
Set-CellXml $t.Cell(4,3) @("25 x 71", "  7    1", "  ----", "2|    |", "5|    |")

Set-CellXml $t.Cell(5,1) @("48 x 70", "  7    0", "  ----", "4|    |", "8|    |")
Set-CellXml $t.Cell(5,2) @("25 x 58", "  5    8", "  ----", "2|    |", "5|    |")
Set-CellXml $t.Cell(5,3) @("40 x 92", "  9    2", "  ----", "4|    |", "0|    |")
